$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 10 with the new contribution entry
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Ryan Conyac"
$ws.Range("D10").Value = "Edited video for sprint 1"

# Update the active selection to D11
$ws.Range("D11").Select()
